$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set G2 to TRUE (was FALSE)
$ws.Range("G2").Value = $true

# Delete entire row 3 (shift cells up) - removes the duplicate "Qty Complete" row
$ws.Rows(3).Delete()

# Update selection to match target diff: entire row 2 selected (activeCell G2, sqref A2:XFD2)
$ws.Rows(2).Select()
